# Apply crypto list updates (prices, volumes, and row40/41 content swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.938.55"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "3.096.37"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.23%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.092.33"
$ws.Range("E8").Value = "  +4.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.59%  "
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.01%  "
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "3.610.15"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").Value = "66.954.10"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("E18").Value = "  +4.02%  "
$ws.Range("D19").Value = "3.096.72"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +18.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +7.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.21%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.05%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "390.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "2.760.93"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.04%  "
